# playerCount in /balance & max player purchase chk
#
# 1. Adds a new AUCTION error code (709) for "max player purchase count
#    reached", reusing the formatting of the row above it.
# 2. Widens the AUCTION description column to fit the new (longer) text.
# 3. Makes AUCTION the active sheet/tab (it was STAT before).

$wb = $excel.ActiveWorkbook

$auction = $wb.Worksheets.Item("AUCTION")
$stat = $wb.Worksheets.Item("STAT")

# --- Add the new error row to AUCTION ------------------------------------
$newRow = 10
$auction.Cells.Item($newRow, 1).Value = 709
$auction.Cells.Item($newRow, 2).Value = "Max player purchase count reached. Cannot buy additional player."

# Match the existing table formatting (bordered/centered code column,
# bordered description column) by copying the format of the row above.
$auction.Range("A9:B9").Copy() | Out-Null
$auction.Range("A10:B10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Widen column B to fit the longer description text.
$auction.Columns.Item(2).ColumnWidth = 60.02

# --- Switch the active sheet / selection from STAT to AUCTION ------------
# (STAT keeps its own selection at B7 but is no longer the active tab, so
# set its selection first, then activate + select on AUCTION last so
# AUCTION ends up as the active/selected tab.)
$stat.Activate() | Out-Null
$stat.Range("B7").Select() | Out-Null

$auction.Activate() | Out-Null
$auction.Range("C13").Select() | Out-Null
